$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    @(" Dubai (DSC)", " October 14 2020", "Capitals won by 13 runs", "Delhi Capitals", "Rajasthan Royals", "Marcus Stoinis ", "18", "19", "1", "0", "94.73"),
    @(" Abu Dhabi", " October 11 2020", "Mumbai won by 5 wickets (with 2 balls remaining)", "Delhi Capitals", "Mumbai Indians", "Marcus Stoinis ", "13", "8", "2", "0", "162.50"),
    @(" Abu Dhabi", " November 02 2020", "Capitals won by 6 wickets (with 6 balls remaining)", "Delhi Capitals", "Royal Challengers Bangalore", "Marcus Stoinis ", "10", "5", "0", "1", "200.00"),
    @(" Dubai (DSC)", " November 05 2020", "Mumbai won by 57 runs", "Delhi Capitals", "Mumbai Indians", "Marcus Stoinis ", "65", "46", "6", "3", "141.30"),
    @(" Dubai (DSC)", " October 27 2020", "Sunrisers won by 88 runs", "Delhi Capitals", "Sunrisers Hyderabad", "Marcus Stoinis ", "5", "6", "1", "0", "83.33"),
    @(" Abu Dhabi", " September 29 2020", "Sunrisers won by 15 runs", "Delhi Capitals", "Sunrisers Hyderabad", "Marcus Stoinis ", "11", "9", "1", "0", "122.22"),
    @(" Abu Dhabi", " October 24 2020", "KKR won by 59 runs", "Delhi Capitals", "Kolkata Knight Riders", "Marcus Stoinis ", "6", "6", "0", "0", "100.00"),
    @(" Sharjah", " October 17 2020", "Capitals won by 5 wickets (with 1 ball remaining)", "Delhi Capitals", "Chennai Super Kings", "Marcus Stoinis ", "24", "14", "1", "2", "171.42"),
    @(" Dubai (DSC)", " October 31 2020", "Mumbai won by 9 wickets (with 34 balls remaining)", "Delhi Capitals", "Mumbai Indians", "Marcus Stoinis ", "2", "3", "0", "0", "66.66"),
    @(" Dubai (DSC)", " November 10 2020", "Mumbai won by 5 wickets (with 8 balls remaining)", "Delhi Capitals", "Mumbai Indians", "Marcus Stoinis ", "0", "1", "0", "0", "0.00"),
    @(" Dubai (DSC)", " October 20 2020", "Kings XI won by 5 wickets (with 6 balls remaining)", "Delhi Capitals", "Kings XI Punjab", "Marcus Stoinis ", "9", "10", "0", "0", "90.00"),
    @(" Dubai (DSC)", " October 05 2020", "Capitals won by 59 runs", "Delhi Capitals", "Royal Challengers Bangalore", "Marcus Stoinis ", "53", "26", "6", "2", "203.84"),
    @(" Dubai (DSC)", " September 20 2020", "Match tied (Capitals won the one-over eliminator)", "Delhi Capitals", "Kings XI Punjab", "Marcus Stoinis ", "53", "21", "7", "3", "252.38"),
    @(" Abu Dhabi", " November 08 2020", "Capitals won by 17 runs", "Delhi Capitals", "Sunrisers Hyderabad", "Marcus Stoinis ", "38", "27", "5", "1", "140.74"),
    @(" Sharjah", " October 03 2020", "Capitals won by 18 runs", "Delhi Capitals", "Kolkata Knight Riders", "Marcus Stoinis ", "1", "3", "0", "0", "33.33"),
    @(" Dubai (DSC)", " September 25 2020", "Capitals won by 44 runs", "Delhi Capitals", "Chennai Super Kings", "Marcus Stoinis ", "5", "3", "1", "0", "166.66"),
    @(" Sharjah", " October 09 2020", "Capitals won by 46 runs", "Delhi Capitals", "Rajasthan Royals", "Marcus Stoinis ", "39", "30", "0", "4", "130.00")
)

$startRow = 19
$ws.Range("G19:K35").NumberFormat = "@"

for ($i = 0; $i -lt $data.Count; $i++) {
    $r = $startRow + $i
    $rowVals = $data[$i]
    for ($c = 0; $c -lt $rowVals.Count; $c++) {
        $ws.Cells.Item($r, $c + 1).Value = $rowVals[$c]
    }
}
